# Updated cryptos list: refresh Price (column D) and Volume(1h) (column E)
# values for each coin row, matching the upstream GitHub Actions data pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "26.883.21"
$ws.Range("E2").Value = "  +1.11%  "
$ws.Range("D3").Value = "1.842.26"
$ws.Range("E3").Value = "  +1.19%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  -0.55%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "309.12"
$ws.Range("E5").Value = "  +1.20%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.005"
$ws.Range("E6").Value = "  -0.34%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4744"
$ws.Range("E7").Value = "  +1.49%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3676"
$ws.Range("E8").Value = "  +2.43%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07201"
$ws.Range("E9").Value = "  +1.08%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9258"
$ws.Range("E10").Value = "  +2.67%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "19.62"
$ws.Range("E11").Value = "  +1.24%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07639"
$ws.Range("E12").Value = "  -2.07%  "
$ws.Range("D13").Value = "1.864.12"
$ws.Range("E13").Value = "  +2.41%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.310"
$ws.Range("E14").Value = "  +1.13%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.407"
$ws.Range("E15").Value = "  +1.22%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "88.73"
$ws.Range("E16").Value = "  +1.69%  "
$ws.Range("E17").Value = "  -0.36%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008641"
$ws.Range("E18").Value = "  +1.11%  "
$ws.Range("E19").Value = "  -0.37%  "
$ws.Range("D20").Value = "26.907.73"
$ws.Range("E20").Value = "  +1.00%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.53"
$ws.Range("E21").Value = "  +2.72%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.044"
$ws.Range("E22").Value = "  +0.80%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.65"
$ws.Range("E23").Value = "  +0.99%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.920"
$ws.Range("E24").Value = "  -1.10%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "152.13"
$ws.Range("E26").Value = "  +1.40%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.001"
$ws.Range("E27").Value = "  +1.72%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "114.25"
$ws.Range("E28").Value = "  +0.69%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.947"
$ws.Range("E29").Value = "  +3.40%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.08850"
$ws.Range("E30").Value = "  +0.57%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.298"
$ws.Range("E31").Value = "  +4.78%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7503"
$ws.Range("E32").Value = "  +3.16%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.177"
$ws.Range("E33").Value = "  +4.84%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.487"
$ws.Range("E34").Value = "  +1.15%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.750"
$ws.Range("E35").Value = "  -0.12%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.089"
$ws.Range("E36").Value = "  +1.44%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.05264"
$ws.Range("E37").Value = "  +3.10%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01948"
$ws.Range("E38").Value = "  +1.33%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.960"
$ws.Range("E39").Value = "  +1.49%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5226"
$ws.Range("E40").Value = "  +3.94%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.932"
$ws.Range("E41").Value = "  +1.44%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1513"
$ws.Range("E42").Value = "  +1.34%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.223"
$ws.Range("E43").Value = "  +3.09%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.56"
$ws.Range("E44").Value = "  +5.52%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4729"
$ws.Range("E45").Value = "  +1.74%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.006"
$ws.Range("E46").Value = "  -0.30%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "101.69"
$ws.Range("E47").Value = "  +3.15%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.606"
$ws.Range("E48").Value = "  +3.34%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "65.55"
$ws.Range("E49").Value = "  +3.01%  "
$ws.Range("E50").Value = "  +0.37%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.8850"
$ws.Range("E51").Value = "  +4.22%  "
